$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new expense row (row 5) below the existing data
# Use Copy for the text-typed columns so Excel doesn't auto-convert the
# date-looking string into a real date serial number.
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("B5").Value = 10
$ws.Range("C3").Copy($ws.Range("C5"))
$ws.Range("D5").Value = "Utilities"
$ws.Range("E5").Value = 1000
